# Update countries & provincias Spain
# Applies the "Datos actualizados" timestamp bump, three country-name
# reorderings in column A (which carry their per-country stats along
# with them), and the refreshed COVID statistics for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 12:13"

# --- Row data: row number -> Country, TotalCases, NewCases, Active, Recovered, Critical, DeathsToday, Deaths ---
$rows = @(
    @{ Row = 4;   Country = "Estados Unidos";         B = 8459041; C = 2388; D = 5504464; E = 2729336; F = 0; G = 19; H = 225241 },
    @{ Row = 20;  Country = "Banglades";               B = 391586;  C = 1380; D = 307141;  E = 78746;   F = 0; G = 18; H = 5699 },
    @{ Row = 42;  Country = "Emiratos Arabes Unidos";  B = 117594;  C = 1077; D = 110313;  E = 6811;    F = 0; G = 4;  H = 470 },
    @{ Row = 43;  Country = "Kuwait";                  B = 116832;  C = 0;    D = 108606;  E = 7516;    F = 0; G = 0;  H = 710 },
    @{ Row = 44;  Country = "Oman";                    B = 111033;  C = 439;  D = 96949;   E = 12962;   F = 0; G = 8;  H = 1122 },
    @{ Row = 59;  Country = "Austria";                 B = 67451;   C = 1524; D = 51407;   E = 15130;   F = 0; G = 10; H = 914 },
    @{ Row = 60;  Country = "Moldavia";                B = 67302;   C = 0;    D = 48493;   E = 17209;   F = 0; G = 0;  H = 1600 },
    @{ Row = 61;  Country = "Armenia";                 B = 66694;   C = 1234; D = 48734;   E = 16859;   F = 0; G = 10; H = 1101 },
    @{ Row = 86;  Country = "Australia";                B = 27429;   C = 30;   D = 25126;   E = 1398;    F = 0; G = 0;  H = 905 },
    @{ Row = 91;  Country = "Malasia";                 B = 22225;   C = 862;  D = 14351;   E = 7681;    F = 0; G = 3;  H = 193 },
    @{ Row = 92;  Country = "Camerun";                 B = 21506;   C = 0;    D = 20117;   E = 965;     F = 0; G = 0;  H = 424 },
    @{ Row = 143; Country = "Estonia";                 B = 4127;    C = 42;   D = 3270;    E = 786;     F = 0; G = 3;  H = 71 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Country
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
}
